# Updated cryptos list on Fri Aug  9 09:12:30 UTC 2024 with GitHub Actions
# Applies the latest crypto market snapshot (prices / 1h volume % changes)
# to the existing "cryptos" worksheet, including the reshuffled ranking
# rows (WrappedeETH dropped out, Kaspa/others shifted up one slot, and
# dogwifhat newly appears at the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume% (E) columns store plain-number-looking values as
# text (e.g. "1.00", "61.162.93"). Force the Text number format first so
# Excel does not silently coerce these assignments into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '61.162.93'
$ws.Range('E2').Value = '  +7.64%  '
$ws.Range('D3').Value = '2.675.15'
$ws.Range('E3').Value = '  +11.45%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '510.62'
$ws.Range('E5').Value = '  +4.73%  '
$ws.Range('D6').Value = '159.41'
$ws.Range('E6').Value = '  +3.47%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  -2.03%  '
$ws.Range('D9').Value = '2.669.85'
$ws.Range('E9').Value = '  +10.32%  '
$ws.Range('D10').Value = '6.50'
$ws.Range('E10').Value = '  +3.54%  '
$ws.Range('E11').Value = '  +6.33%  '
$ws.Range('E12').Value = '  +4.45%  '
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('D14').Value = '3.131.82'
$ws.Range('E14').Value = '  +10.73%  '
$ws.Range('D15').Value = '61.119.47'
$ws.Range('E15').Value = '  +7.08%  '
$ws.Range('D16').Value = '21.90'
$ws.Range('E16').Value = '  +6.70%  '
$ws.Range('E17').Value = '  +7.12%  '
$ws.Range('D18').Value = '2.669.34'
$ws.Range('E18').Value = '  +10.14%  '
$ws.Range('D19').Value = '4.83'
$ws.Range('E19').Value = '  +2.54%  '
$ws.Range('D20').Value = '349.44'
$ws.Range('E20').Value = '  +7.95%  '
$ws.Range('D21').Value = '10.60'
$ws.Range('E21').Value = '  +6.69%  '
$ws.Range('D22').Value = '6.23'
$ws.Range('E22').Value = '  +5.57%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = '60.89'
$ws.Range('E24').Value = '  +5.30%  '
$ws.Range('D25').Value = '0.427'
$ws.Range('E25').Value = '  +5.97%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.168'
$ws.Range('E26').Value = '  +4.29%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0877'
$ws.Range('E28').Value = '  +12.95%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '7.65'
$ws.Range('E29').Value = '  +5.86%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '19.62'
$ws.Range('E31').Value = '  +5.98%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '156.71'
$ws.Range('E32').Value = '  +4.44%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.58'
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = '5.87'
$ws.Range('E34').Value = '  +11.47%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.09'
$ws.Range('E35').Value = '  +9.37%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.23'
$ws.Range('E36').Value = '  +6.90%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = '317.64'
$ws.Range('E37').Value = '  +16.49%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = '0.866'
$ws.Range('E38').Value = '  +3.97%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.51'
$ws.Range('E39').Value = '  +10.59%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.80'
$ws.Range('E40').Value = '  +7.91%  '
$ws.Range('E41').Value = '  +30.79%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '35.45'
$ws.Range('E42').Value = '  +4.15%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '0.0582'
$ws.Range('E43').Value = '  +10.19%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.639'
$ws.Range('E44').Value = '  +7.19%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = '0.101'
$ws.Range('E45').Value = '  -1.56%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '19.77'
$ws.Range('E47').Value = '  +13.48%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '4.93'
$ws.Range('E48').Value = '  +10.56%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.075.28'
$ws.Range('E49').Value = '  +10.19%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0237'
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '1.88'
$ws.Range('E51').Value = '  +12.61%  '
